$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused template rows (10-12), keeping only the single
# data row (row 9) under the header.
$ws.Rows("10:12").Delete()

# Row 9 becomes a "blank template" row: no name, new phone number, new
# attachment path.
$ws.Range("A9").Value = ""
$ws.Range("B9").Value = 31994773182
$ws.Range("C9").Value = "C:\Projetos\bot whatsapp\test.jpg"

$wb.Save()
